# Apply the "Updated symbol list" data refresh (coinranking.com scrape, 9 Jan 2023 19:xx UTC)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '278.31'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '4.71%'

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '26.87'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '0.50%'

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '4.914'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '4.78%'

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.06392'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '4.97%'

# Row 6
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '3.94%'

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '3.356'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '5.91%'

# Row 8
$ws.Range("B8").Value = 'FTXToken'
$ws.Range("C8").Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '1.210'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '33.63%'

# Row 9
$ws.Range("B9").Value = 'MXToken'
$ws.Range("C9").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.8859'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '4.25%'

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.1483'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '5.12%'

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.05271'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '4.98%'

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.07420'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '4.72%'

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.03134'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '-1.33%'

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.09062'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '0.46%'

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.001567'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '2.56%'

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.0006345'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '4.50%'

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.006015'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '0.39%'

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.489'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '0.81%'

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '2.281'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '0.77%'

# Row 20
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '2.20%'

# Row 21
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '2.48%'

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '3.915'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '-4.10%'

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.04342'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '2.57%'

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.001180'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '-0.39%'

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.003677'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '-11.12%'

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.0001200'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '-0.01%'

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.0001616'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '-3.89%'

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.04072'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '4.09%'

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.006655'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '59.10%'

# Row 42
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '5.34%'

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.002361'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '11.30%'

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.01281'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '1.65%'

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.00005268'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '3.15%'

# Row 47
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '1,450.06%'

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.02121'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '-13.35%'

# Row 50
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '-0.12%'
